$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header text: Score -> Score E1, Notes -> Notes E1
$ws.Range("D1").Value = "Score E1"
$ws.Range("E1").Value = "Notes E1"

# Update the view: scroll so column C is the leftmost visible column, and select D2
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("D2").Select()
